# Refresh the cryptocurrency price/volume snapshot (coinranking.com scrape).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '62.437.01'
$ws.Cells.Item(2, 5).Value = '  -2.00%  '

$ws.Cells.Item(3, 4).Value = '3.172.84'
$ws.Cells.Item(3, 5).Value = '  -3.58%  '

$ws.Cells.Item(4, 5).Value = '  +0.05%  '

$ws.Cells.Item(5, 4).Value = '''587.09'
$ws.Cells.Item(5, 5).Value = '  -2.92%  '

$ws.Cells.Item(6, 4).Value = '''135.07'
$ws.Cells.Item(6, 5).Value = '  -4.87%  '

$ws.Cells.Item(7, 5).Value = '  -0.02%  '

$ws.Cells.Item(8, 4).Value = '3.170.22'
$ws.Cells.Item(8, 5).Value = '  -3.54%  '

$ws.Cells.Item(9, 5).Value = '  -2.17%  '

$ws.Cells.Item(10, 5).Value = '  -5.36%  '

$ws.Cells.Item(11, 5).Value = '  -2.91%  '

$ws.Cells.Item(12, 4).Value = '''0.454'
$ws.Cells.Item(12, 5).Value = '  -3.40%  '

$ws.Cells.Item(13, 4).Value = '''0.0000234'
$ws.Cells.Item(13, 5).Value = '  -5.05%  '

$ws.Cells.Item(14, 4).Value = '''33.25'
$ws.Cells.Item(14, 5).Value = '  -3.79%  '

$ws.Cells.Item(15, 4).Value = '3.695.78'
$ws.Cells.Item(15, 5).Value = '  -3.65%  '

$ws.Cells.Item(16, 5).Value = '  -1.75%  '

$ws.Cells.Item(17, 4).Value = '3.175.99'
$ws.Cells.Item(17, 5).Value = '  -3.56%  '

$ws.Cells.Item(18, 4).Value = '62.466.21'
$ws.Cells.Item(18, 5).Value = '  -2.06%  '

$ws.Cells.Item(19, 4).Value = '''6.51'
$ws.Cells.Item(19, 5).Value = '  -4.81%  '

$ws.Cells.Item(20, 4).Value = '''455.51'
$ws.Cells.Item(20, 5).Value = '  -5.12%  '

$ws.Cells.Item(21, 5).Value = '  -1.10%  '

$ws.Cells.Item(22, 4).Value = '''0.701'
$ws.Cells.Item(22, 5).Value = '  -3.81%  '

$ws.Cells.Item(23, 5).Value = '  -4.93%  '

$ws.Cells.Item(24, 4).Value = '''83.52'
$ws.Cells.Item(24, 5).Value = '  -0.99%  '

$ws.Cells.Item(25, 4).Value = '''13.24'
$ws.Cells.Item(25, 5).Value = '  -1.85%  '

$ws.Cells.Item(26, 5).Value = '  +0.05%  '

$ws.Cells.Item(27, 4).Value = '''2.68'
$ws.Cells.Item(27, 5).Value = '  -3.33%  '

$ws.Cells.Item(28, 5).Value = '  +0.03%  '

$ws.Cells.Item(29, 4).Value = '''6.84'
$ws.Cells.Item(29, 5).Value = '  -6.10%  '

$ws.Cells.Item(30, 4).Value = '''7.73'

$ws.Cells.Item(31, 4).Value = '''2.01'
$ws.Cells.Item(31, 5).Value = '  -7.20%  '

$ws.Cells.Item(32, 4).Value = '''27.25'
$ws.Cells.Item(32, 5).Value = '  -5.77%  '

$ws.Cells.Item(33, 5).Value = '  -2.45%  '

$ws.Cells.Item(34, 5).Value = '  -6.13%  '

$ws.Cells.Item(35, 5).Value = '  -6.22%  '

$ws.Cells.Item(36, 5).Value = '  -1.11%  '

$ws.Cells.Item(37, 4).Value = '''51.10'
$ws.Cells.Item(37, 5).Value = '  -3.63%  '

$ws.Cells.Item(38, 4).Value = '0.0₃0697'
$ws.Cells.Item(38, 5).Value = '  -6.64%  '

$ws.Cells.Item(39, 5).Value = '  -3.76%  '

$ws.Cells.Item(40, 4).Value = '''2.72'
$ws.Cells.Item(40, 5).Value = '  -0.26%  '

$ws.Cells.Item(41, 5).Value = '  +0.73%  '

$ws.Cells.Item(42, 4).Value = '''7.98'
$ws.Cells.Item(42, 5).Value = '  -4.42%  '

$ws.Cells.Item(43, 4).Value = '''390.31'
$ws.Cells.Item(43, 5).Value = '  -8.19%  '

$ws.Cells.Item(44, 4).Value = '2.797.58'
$ws.Cells.Item(44, 5).Value = '  -8.30%  '

$ws.Cells.Item(45, 2).Value = 'TheGraph'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(45, 4).Value = '''0.249'
$ws.Cells.Item(45, 5).Value = '  -5.79%  '

$ws.Cells.Item(46, 2).Value = 'USDe'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(46, 4).Value = '''0.999'
$ws.Cells.Item(46, 5).Value = '  -0.05%  '

$ws.Cells.Item(47, 5).Value = '  -2.57%  '

$ws.Cells.Item(48, 2).Value = 'Arweave'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(48, 4).Value = '''35.75'
$ws.Cells.Item(48, 5).Value = '  +3.04%  '

$ws.Cells.Item(49, 4).Value = '''124.80'
$ws.Cells.Item(49, 5).Value = '  +0.26%  '

$ws.Cells.Item(50, 4).Value = '''25.24'
$ws.Cells.Item(50, 5).Value = '  -3.64%  '

$ws.Cells.Item(51, 5).Value = '  -3.77%  '
